$d = $word.ActiveDocument

# Locate the already-empty paragraph (no text, no inline image) that must be preserved.
$keepPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $r = $p.Range
    if ($r.Text -eq [char]13 -and $r.InlineShapes.Count -eq 0) {
        $keepPara = $p
        break
    }
}

$startKeep = $keepPara.Range.Start
$endKeep = $keepPara.Range.End

# Remove all content after the paragraph we keep.
if ($endKeep -lt $d.Content.End) {
    $afterRange = $d.Range($endKeep, $d.Content.End)
    $afterRange.Delete()
}

# Remove all content before the paragraph we keep.
if ($startKeep -gt 0) {
    $beforeRange = $d.Range(0, $startKeep)
    $beforeRange.Delete()
}

Write-Host "Done. Paragraph count now:" $d.Paragraphs.Count
